# Fruta / hortaliza, semanal
# Insert a new weekly record at row 48 of the "Berenjena" data sheet,
# pushing the existing rows 48-65 down to 49-66 (dimension grows to A1:R66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 48 - this shifts rows 48..65
# down to 49..66, preserving all of their existing data/formatting.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly observation.
$ws.Cells.Item(48, 1).Value  = 7
$ws.Cells.Item(48, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(48, 3).Value  = "Ñuble"
$ws.Cells.Item(48, 4).Value  = 45029
$ws.Cells.Item(48, 5).Value  = 16
$ws.Cells.Item(48, 6).Value  = 100112001
$ws.Cells.Item(48, 7).Value  = "Berenjena"
$ws.Cells.Item(48, 8).Value  = "Sin especificar"
$ws.Cells.Item(48, 9).Value  = "Primera"
$ws.Cells.Item(48, 10).Value = 50
$ws.Cells.Item(48, 11).Value = 10000
$ws.Cells.Item(48, 12).Value = 10000
$ws.Cells.Item(48, 13).Value = 10000
$ws.Cells.Item(48, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 167
$ws.Cells.Item(48, 17).Value = 60
$ws.Cells.Item(48, 18).Value = "Hortaliza"
